# Scheduled runner update: refresh computed market-price / profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ -- columns H:N) on the affected leve rows across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets. Cells with no value in the
# refreshed data (e.g. a profit column that no longer applies) are cleared
# via $null so the cell is dropped rather than left as a stale 0.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 274.3889
$ws.Range("I33").Value = 149.11111
$ws.Range("J33").Value = 399.66666
$ws.Range("K33").Value = 149.11111
$ws.Range("L33").Value = 399.66666
$ws.Range("M33").Value = 79.88889
$ws.Range("N33").Value = -857.66666
$ws.Range("H86").Value = 100004380
$ws.Range("I86").Value = 3406.8572
$ws.Range("K86").Value = 3406.8572
$ws.Range("M86").Value = -2283.8572
$ws.Range("H89").Value = 100004380
$ws.Range("I89").Value = 3406.8572
$ws.Range("K89").Value = 17034.286
$ws.Range("M89").Value = -11418.286
$ws.Range("H111").Value = 1051.5
$ws.Range("I111").Value = 1269.25
$ws.Range("J111").Value = 616
$ws.Range("K111").Value = 3807.75
$ws.Range("L111").Value = 1848
$ws.Range("M111").Value = -740.75
$ws.Range("N111").Value = -7982
$ws.Range("H121").Value = 772.9231
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 768
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 2304
$ws.Range("M121").Value = -653
$ws.Range("N121").Value = -5798
$ws.Range("H125").Value = 817.7778
$ws.Range("I125").Value = 832.8570999999999
$ws.Range("J125").Value = 765
$ws.Range("K125").Value = 7495.7139
$ws.Range("L125").Value = 6885
$ws.Range("M125").Value = -5035.7139
$ws.Range("N125").Value = -11805
$ws.Range("H135").Value = 2342.55
$ws.Range("I135").Value = 1602.8334
$ws.Range("J135").Value = 9000
$ws.Range("K135").Value = 14425.5006
$ws.Range("L135").Value = 81000
$ws.Range("M135").Value = -11890.5006
$ws.Range("N135").Value = -86070
$ws.Range("H137").Value = 1258.8334
$ws.Range("I137").Value = 825
$ws.Range("K137").Value = 2475
$ws.Range("M137").Value = 75
$ws.Range("H138").Value = 2754.83
$ws.Range("I138").Value = 951.5172
$ws.Range("J138").Value = 3491.3943
$ws.Range("K138").Value = 2854.5516
$ws.Range("L138").Value = 10474.1829
$ws.Range("M138").Value = 2285.4484
$ws.Range("N138").Value = -20754.1829

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1374.7428
$ws.Range("J74").Value = 1850.25
$ws.Range("L74").Value = 1850.25
$ws.Range("N74").Value = -3598.25
$ws.Range("H77").Value = 1374.7428
$ws.Range("J77").Value = 1850.25
$ws.Range("L77").Value = 9251.25
$ws.Range("N77").Value = -17987.25
$ws.Range("H102").Value = 3922
$ws.Range("I102").Value = 4277.5
$ws.Range("K102").Value = 4277.5
$ws.Range("M102").Value = -2655.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 6426.6665
$ws.Range("J74").Value = 6426.6665
$ws.Range("L74").Value = 6426.6665
$ws.Range("N74").Value = -8298.666499999999
$ws.Range("H77").Value = 6426.6665
$ws.Range("J77").Value = 6426.6665
$ws.Range("L77").Value = 19279.9995
$ws.Range("N77").Value = -28639.9995
$ws.Range("H105").Value = 8308.25
$ws.Range("I105").Value = 6336.273
$ws.Range("J105").Value = 30000
$ws.Range("K105").Value = 6336.273
$ws.Range("L105").Value = 30000
$ws.Range("M105").Value = -4589.273
$ws.Range("N105").Value = -33494
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = $null
$ws.Range("N133").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 3298
$ws.Range("I81").Value = 3298
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3298
$ws.Range("L81").Value = $null
$ws.Range("N81").Value = 0
$ws.Range("M81").Value = -2300
$ws.Range("H84").Value = 3298
$ws.Range("I84").Value = 3298
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9894
$ws.Range("L84").Value = $null
$ws.Range("N84").Value = 0
$ws.Range("M84").Value = -4902
$ws.Range("H105").Value = 2754.7646
$ws.Range("I105").Value = 2621
$ws.Range("K105").Value = 2621
$ws.Range("M105").Value = -874

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 109.77778
$ws.Range("I6").Value = 109.77778
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 329.33334
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = -216.33334
$ws.Range("H110").Value = 13874.5
$ws.Range("I110").Value = 4996
$ws.Range("J110").Value = 15142.857
$ws.Range("K110").Value = 14988
$ws.Range("L110").Value = 45428.571
$ws.Range("M110").Value = -10898
$ws.Range("N110").Value = -53608.571
$ws.Range("H113").Value = 30043.266
$ws.Range("I113").Value = 870.25
$ws.Range("J113").Value = 39019.58
$ws.Range("K113").Value = 2610.75
$ws.Range("L113").Value = 117058.74
$ws.Range("M113").Value = -440.75
$ws.Range("N113").Value = -121398.74
$ws.Range("H122").Value = 1005.28125
$ws.Range("J122").Value = 1022.3
$ws.Range("L122").Value = 9200.699999999999
$ws.Range("N122").Value = -14100.7
$ws.Range("H132").Value = 1069.4375
$ws.Range("I132").Value = 884.1818
$ws.Range("J132").Value = 1477
$ws.Range("K132").Value = 7957.6362
$ws.Range("L132").Value = 13293
$ws.Range("M132").Value = -5427.6362
$ws.Range("N132").Value = -18353

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3637
$ws.Range("I80").Value = 3451.0588
$ws.Range("J80").Value = 3988.2222
$ws.Range("K80").Value = 3451.0588
$ws.Range("L80").Value = 3988.2222
$ws.Range("M80").Value = -2453.0588
$ws.Range("N80").Value = -5984.2222
$ws.Range("H83").Value = 3637
$ws.Range("I83").Value = 3451.0588
$ws.Range("J83").Value = 3988.2222
$ws.Range("K83").Value = 17255.294
$ws.Range("L83").Value = 19941.111
$ws.Range("M83").Value = -12263.294
$ws.Range("N83").Value = -29925.111
$ws.Range("H108").Value = 35000
$ws.Range("J108").Value = 35000
$ws.Range("L108").Value = 35000
$ws.Range("N108").Value = -42680
$ws.Range("H126").Value = 4817.7144
$ws.Range("I126").Value = 4817.7144
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14453.1432
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = -11983.1432
$ws.Range("H132").Value = 3912.4285
$ws.Range("I132").Value = 3734.375
$ws.Range("J132").Value = 4482.2
$ws.Range("K132").Value = 11203.125
$ws.Range("L132").Value = 13446.6
$ws.Range("M132").Value = -8673.125
$ws.Range("N132").Value = -18506.6
$ws.Range("H136").Value = 9437.684999999999
$ws.Range("J136").Value = 9437.684999999999
$ws.Range("L136").Value = 28313.055
$ws.Range("N136").Value = -33413.055
$ws.Range("H140").Value = 56120
$ws.Range("J140").Value = 56120
$ws.Range("L140").Value = 56120
$ws.Range("N140").Value = -66480

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1998.8
$ws.Range("I7").Value = 2012.6818
$ws.Range("J7").Value = 1960.625
$ws.Range("K7").Value = 2012.6818
$ws.Range("L7").Value = 1960.625
$ws.Range("M7").Value = -1900.6818
$ws.Range("N7").Value = -2184.625
$ws.Range("H40").Value = 562262.2
$ws.Range("I40").Value = 842626.2
$ws.Range("K40").Value = 842626.2
$ws.Range("M40").Value = -842490.2
$ws.Range("H82").Value = 1698.0555
$ws.Range("I82").Value = 1674.2354
$ws.Range("J82").Value = 2103
$ws.Range("K82").Value = 1674.2354
$ws.Range("L82").Value = 2103
$ws.Range("M82").Value = -1313.2354
$ws.Range("N82").Value = -2825
$ws.Range("H85").Value = 1698.0555
$ws.Range("I85").Value = 1674.2354
$ws.Range("J85").Value = 2103
$ws.Range("K85").Value = 1674.2354
$ws.Range("L85").Value = 2103
$ws.Range("M85").Value = -426.2354
$ws.Range("N85").Value = -4599
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -10900
$ws.Range("H126").Value = 1998.8
$ws.Range("I126").Value = 2012.6818
$ws.Range("J126").Value = 1960.625
$ws.Range("K126").Value = 6038.0454
$ws.Range("L126").Value = 5881.875
$ws.Range("M126").Value = -3568.0454
$ws.Range("N126").Value = -10821.875
$ws.Range("H132").Value = 1995.7931
$ws.Range("I132").Value = 1309.5264
$ws.Range("J132").Value = 3299.7
$ws.Range("K132").Value = 3928.5792
$ws.Range("L132").Value = 9899.099999999999
$ws.Range("M132").Value = -1398.5792
$ws.Range("N132").Value = -14959.1
$ws.Range("H136").Value = 4508.815
$ws.Range("I136").Value = 2076.9
$ws.Range("K136").Value = 6230.700000000001
$ws.Range("M136").Value = -3680.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 49000
$ws.Range("J18").Value = 49000
$ws.Range("L18").Value = 49000
$ws.Range("N18").Value = -49346
$ws.Range("H42").Value = 54524.5
$ws.Range("J42").Value = 54524.5
$ws.Range("L42").Value = 54524.5
$ws.Range("N42").Value = -55280.5
$ws.Range("H43").Value = 19250
$ws.Range("J43").Value = 14000
$ws.Range("L43").Value = 14000
$ws.Range("N43").Value = -14298
